$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume cells keep their original text formatting
# (these are stored as text, e.g. "301.97" / "4.99%", not numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "301.97"
$ws.Range("E2").Value = "4.99%"
$ws.Range("D3").Value = "34.83"
$ws.Range("E3").Value = "12.44%"
$ws.Range("D4").Value = "5.181"
$ws.Range("E4").Value = "5.36%"
$ws.Range("D5").Value = "0.07761"
$ws.Range("E5").Value = "6.04%"
$ws.Range("D6").Value = "2.277"
$ws.Range("E6").Value = "-3.31%"
$ws.Range("D7").Value = "8.023"
$ws.Range("E7").Value = "3.83%"
$ws.Range("D8").Value = "4.004"
$ws.Range("E8").Value = "7.51%"
$ws.Range("D9").Value = "0.9281"
$ws.Range("E9").Value = "2.81%"
$ws.Range("D10").Value = "0.1011"
$ws.Range("E10").Value = "10.83%"
$ws.Range("D11").Value = "0.1816"
$ws.Range("E11").Value = "7.38%"
$ws.Range("D12").Value = "0.08536"
$ws.Range("E12").Value = "4.70%"
$ws.Range("D13").Value = "0.03457"
$ws.Range("E13").Value = "10.48%"
$ws.Range("D14").Value = "0.09911"
$ws.Range("E14").Value = "-0.27%"
$ws.Range("D15").Value = "0.001484"
$ws.Range("E15").Value = "-0.92%"
$ws.Range("D16").Value = "0.04622"
$ws.Range("E16").Value = "2.57%"
$ws.Range("D17").Value = "0.005751"
$ws.Range("E17").Value = "-0.44%"
$ws.Range("D18").Value = "3.472"
$ws.Range("E18").Value = "-0.64%"
$ws.Range("D19").Value = "2.106"
$ws.Range("D20").Value = "0.3438"
$ws.Range("E20").Value = "3.24%"
$ws.Range("E21").Value = "2.65%"
$ws.Range("D22").Value = "4.611"
$ws.Range("E22").Value = "9.90%"
$ws.Range("D23").Value = "0.2343"
$ws.Range("E23").Value = "11.58%"
$ws.Range("D24").Value = "0.001224"
$ws.Range("E24").Value = "1.12%"
$ws.Range("D25").Value = "0.004410"
$ws.Range("E25").Value = "5.95%"
$ws.Range("D26").Value = "0.0001307"
$ws.Range("E26").Value = "0.50%"
$ws.Range("D27").Value = "0.0003422"
$ws.Range("E27").Value = "0.83%"
$ws.Range("D39").Value = "0.01759"
$ws.Range("E39").Value = "11.62%"
$ws.Range("D40").Value = "0.04709"
$ws.Range("E40").Value = "5.99%"
$ws.Range("D41").Value = "0.007618"
$ws.Range("E41").Value = "3.86%"
$ws.Range("E42").Value = "5.86%"
$ws.Range("D43").Value = "0.006857"
$ws.Range("E43").Value = "-27.97%"
$ws.Range("D44").Value = "0.002213"
$ws.Range("E44").Value = "-0.77%"
$ws.Range("D45").Value = "0.009196"
$ws.Range("E45").Value = "2.17%"
$ws.Range("D46").Value = "0.00005938"
$ws.Range("E46").Value = "-2.71%"
$ws.Range("E47").Value = "0.30%"
$ws.Range("E48").Value = "17.28%"
$ws.Range("D49").Value = "0.002706"
$ws.Range("E49").Value = "35.22%"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").Value = "0.30%"
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").Value = "0.30%"
